$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update N (JKL) and O (OS) columns to False for all "sleep" rows (2025-02-01 .. 2025-02-18)
$sleepRows = @(2, 5, 8, 11, 14, 17, 20, 23, 26, 29, 32, 35, 38, 41, 44, 47, 50, 53)
foreach ($r in $sleepRows) {
    $ws.Cells.Item($r, 14).Value = $false
    $ws.Cells.Item($r, 15).Value = $false
}

# Row 53 (2025-02-18, sleep): StayStrong (I) changes from True to False
$ws.Cells.Item(53, 9).Value = $false

# Row 54 (2025-02-18, activity): StayStrong (I) changes from False to True
$ws.Cells.Item(54, 9).Value = $true

# Append new rows for 2025-02-19 (sleep, activity, weekly_activity)
$ws.Cells.Item(56, 1).Value = "'2025-02-19"
$ws.Cells.Item(56, 1).Style = "Normal"
$ws.Cells.Item(56, 2).Value = "sleep"
$ws.Cells.Item(56, 3).Value = $true
$ws.Cells.Item(56, 4).Value = $false
$ws.Cells.Item(56, 5).Value = $true
$ws.Cells.Item(56, 6).Value = $true
$ws.Cells.Item(56, 7).Value = $true
$ws.Cells.Item(56, 8).Value = $true
$ws.Cells.Item(56, 9).Value = $true
$ws.Cells.Item(56, 10).Value = $false
$ws.Cells.Item(56, 11).Value = $true
$ws.Cells.Item(56, 12).Value = $true
$ws.Cells.Item(56, 13).Value = $true
$ws.Cells.Item(56, 14).Value = $false
$ws.Cells.Item(56, 15).Value = $false

$ws.Cells.Item(57, 1).Value = "'2025-02-19"
$ws.Cells.Item(57, 1).Style = "Normal"
$ws.Cells.Item(57, 2).Value = "activity"
$ws.Cells.Item(57, 3).Value = $true
$ws.Cells.Item(57, 4).Value = $false
$ws.Cells.Item(57, 5).Value = $true
$ws.Cells.Item(57, 6).Value = $false
$ws.Cells.Item(57, 7).Value = $true
$ws.Cells.Item(57, 8).Value = $true
$ws.Cells.Item(57, 9).Value = $true
$ws.Cells.Item(57, 10).Value = $true
$ws.Cells.Item(57, 11).Value = $false
$ws.Cells.Item(57, 12).Value = $false
$ws.Cells.Item(57, 13).Value = $true
$ws.Cells.Item(57, 14).Value = $false
$ws.Cells.Item(57, 15).Value = $false

$ws.Cells.Item(58, 1).Value = "'2025-02-19"
$ws.Cells.Item(58, 1).Style = "Normal"
$ws.Cells.Item(58, 2).Value = "weekly_activity"
$ws.Cells.Item(58, 3).Value = $false
$ws.Cells.Item(58, 4).Value = $false
$ws.Cells.Item(58, 5).Value = $true
$ws.Cells.Item(58, 6).Value = $false
$ws.Cells.Item(58, 7).Value = $true
$ws.Cells.Item(58, 8).Value = $false
$ws.Cells.Item(58, 9).Value = $true
$ws.Cells.Item(58, 10).Value = $true
$ws.Cells.Item(58, 11).Value = $false
$ws.Cells.Item(58, 12).Value = $false
$ws.Cells.Item(58, 13).Value = $true
$ws.Cells.Item(58, 14).Value = $false
$ws.Cells.Item(58, 15).Value = $false
